$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2:L2").NumberFormat = "@"

$ws.Range("A2").Value = "1"
$ws.Range("B2").Value = "Гаджиев Расул Магомедович"
$ws.Range("C2").Value = "asdasdasdasd"
$ws.Range("E2").Value = "1в"
$ws.Range("F2").Value = "Образование"
$ws.Range("I2").Value = "1"
$ws.Range("J2").Value = "12.12.2023"
$ws.Range("K2").Value = "30.03.2023"
$ws.Range("L2").Value = "30.03.2023"
